$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on the Price cells that are changing so
# Excel does not auto-convert numeric-looking strings (e.g. "1.000") into numbers.
$dCells = @("D2","D3","D5","D6","D7","D8","D11","D12","D13","D14","D15","D16","D17","D18","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D33","D34","D35","D36","D37","D38","D40","D41","D42","D43","D44","D45","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) and Volume(1h) (E) values row by row.
$ws.Range("D2").Value = "29.221.24"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.857.75"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "0.7116"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").Value = "237.81"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.08120"
$ws.Range("E8").Value = "  +8.82%  "
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D11").Value = "0.08194"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "1.850.11"
$ws.Range("E12").Value = "  -2.49%  "
$ws.Range("D13").Value = "5.169"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").Value = "0.7073"
$ws.Range("E14").Value = "  -3.01%  "
$ws.Range("D15").Value = "89.41"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "29.225.60"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").Value = "5.791"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "0.000007890"
$ws.Range("E18").Value = "  +2.88%  "
$ws.Range("E19").Value = "  +1.66%  "
$ws.Range("D20").Value = "237.48"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "2.103.85"
$ws.Range("E22").Value = "  -1.31%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "7.412"
$ws.Range("E24").Value = "  -2.69%  "
$ws.Range("D25").Value = "162.49"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D26").Value = "0.1462"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "8.963"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").Value = "18.08"
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("D29").Value = "1.960"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").Value = "1.427"
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("D31").Value = "4.400"
$ws.Range("E31").Value = "  -2.76%  "
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").Value = "4.019"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "0.05218"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "1.167"
$ws.Range("E35").Value = "  -2.06%  "
$ws.Range("D36").Value = "0.7081"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").Value = "0.9990"
$ws.Range("E37").Value = "  -4.28%  "
$ws.Range("D38").Value = "2.675"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("E39").Value = "  -0.63%  "
$ws.Range("D40").Value = "2.728"
$ws.Range("E40").Value = "  +1.90%  "
$ws.Range("D41").Value = "1.142.22"
$ws.Range("E41").Value = "  +7.05%  "
$ws.Range("D42").Value = "0.9226"
$ws.Range("E42").Value = "  -2.14%  "
$ws.Range("D43").Value = "0.4281"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").Value = "5.868"
$ws.Range("E44").Value = "  -3.27%  "
$ws.Range("D45").Value = "70.03"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "102.37"
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("D48").Value = "1.775"
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("D49").Value = "2.008.23"
$ws.Range("E49").Value = "  -2.61%  "
$ws.Range("D50").Value = "9.193"
$ws.Range("E50").Value = "  +0.79%  "
$ws.Range("D51").Value = "6.955"
$ws.Range("E51").Value = "  -1.55%  "

# Restore the cell style to Normal so no stray number-format style lingers
# on cells that should remain unstyled, matching the original workbook look.
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
